# NATMI re-run: a new "MuSCs" target cluster was added to the TPM analysis,
# which changes the specificity denominators for the existing ECs->FAPs and
# FAPs->FAPs rows, and adds two new rows (ECs->MuSCs, FAPs->MuSCs).
# Final layout (rows 2-5), sending/ligand/receptor/target clusters:
#   row2: ECs  / Rln3 / Rxfp2 / FAPs
#   row3: ECs  / Rln3 / Rxfp2 / MuSCs   (new)
#   row4: FAPs / Rln3 / Rxfp2 / FAPs    (was row 3)
#   row5: FAPs / Rln3 / Rxfp2 / MuSCs   (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> FAPs (existing row, only some specificity columns change) ---
$ws.Range("M2").Value = 0.06617233333333333
$ws.Range("O2").Value = 0.9596962108540322
$ws.Range("P2").Value = 0.9596962108540322
$ws.Range("Q2").Value = 0.01489519371633333
$ws.Range("S2").Value = 0.1842797290271686
$ws.Range("T2").Value = 0.1842797290271685

# --- Row 3 (new): ECs -> MuSCs ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rln3"
$ws.Range("C3").Value = "Rxfp2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.225097
$ws.Range("H3").Value = 0.675291
$ws.Range("I3").Value = 0.1920188148530651
$ws.Range("J3").Value = 0.1920188148530651
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002779
$ws.Range("N3").Value = 0.008337000000000001
$ws.Range("O3").Value = 0.04030378914596769
$ws.Range("P3").Value = 0.04030378914596769
$ws.Range("Q3").Value = 0.000625544563
$ws.Range("R3").Value = 0.005629901067
$ws.Range("S3").Value = 0.007739085825896545
$ws.Range("T3").Value = 0.007739085825896544

# --- Row 4 (was row 3): FAPs -> FAPs ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rln3"
$ws.Range("C4").Value = "Rxfp2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9471683333333334
$ws.Range("H4").Value = 2.841505
$ws.Range("I4").Value = 0.807981185146935
$ws.Range("J4").Value = 0.8079811851469348
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.06617233333333333
$ws.Range("N4").Value = 0.198517
$ws.Range("O4").Value = 0.9596962108540322
$ws.Range("P4").Value = 0.9596962108540322
$ws.Range("Q4").Value = 0.06267633867611111
$ws.Range("R4").Value = 0.564087048085
$ws.Range("S4").Value = 0.7754164818268637
$ws.Range("T4").Value = 0.7754164818268636

# --- Row 5 (new): FAPs -> MuSCs ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rln3"
$ws.Range("C5").Value = "Rxfp2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.9471683333333334
$ws.Range("H5").Value = 2.841505
$ws.Range("I5").Value = 0.807981185146935
$ws.Range("J5").Value = 0.8079811851469348
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.002779
$ws.Range("N5").Value = 0.008337000000000001
$ws.Range("O5").Value = 0.04030378914596769
$ws.Range("P5").Value = 0.04030378914596769
$ws.Range("Q5").Value = 0.002632180798333333
$ws.Range("R5").Value = 0.023689627185
$ws.Range("S5").Value = 0.03256470332007114
$ws.Range("T5").Value = 0.03256470332007114
